$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("61+37=", $true, $true, $false, $false, $false, $true, 1, $false, "39-5=", 2)
$null = $d.Content.Find.Execute("89-53=", $true, $true, $false, $false, $false, $true, 1, $false, "37+27=", 2)
$null = $d.Content.Find.Execute("67-35=", $true, $true, $false, $false, $false, $true, 1, $false, "48+2=", 2)
$null = $d.Content.Find.Execute("74+20=", $true, $true, $false, $false, $false, $true, 1, $false, "13+65=", 2)
$null = $d.Content.Find.Execute("36-21=", $true, $true, $false, $false, $false, $true, 1, $false, "21+30=", 2)
$null = $d.Content.Find.Execute("55-1=", $true, $true, $false, $false, $false, $true, 1, $false, "49-3=", 2)
$null = $d.Content.Find.Execute("92-44=", $true, $true, $false, $false, $false, $true, 1, $false, "67+10=", 2)
$null = $d.Content.Find.Execute("51-42=", $true, $true, $false, $false, $false, $true, 1, $false, "50+1=", 2)
$null = $d.Content.Find.Execute("96-70=", $true, $true, $false, $false, $false, $true, 1, $false, "77+4=", 2)
$null = $d.Content.Find.Execute("47+52=", $true, $true, $false, $false, $false, $true, 1, $false, "51+44=", 2)
$null = $d.Content.Find.Execute("42+17=", $true, $true, $false, $false, $false, $true, 1, $false, "76-15=", 2)
$null = $d.Content.Find.Execute("25+44=", $true, $true, $false, $false, $false, $true, 1, $false, "66+0=", 2)
$null = $d.Content.Find.Execute("7+55=", $true, $true, $false, $false, $false, $true, 1, $false, "8+0=", 2)
$null = $d.Content.Find.Execute("56-43=", $true, $true, $false, $false, $false, $true, 1, $false, "69-32=", 2)
$null = $d.Content.Find.Execute("82-12=", $true, $true, $false, $false, $false, $true, 1, $false, "8+59=", 2)
$null = $d.Content.Find.Execute("24+73=", $true, $true, $false, $false, $false, $true, 1, $false, "96-48=", 2)
$null = $d.Content.Find.Execute("63-8=", $true, $true, $false, $false, $false, $true, 1, $false, "65+12=", 2)
$null = $d.Content.Find.Execute("44+7=", $true, $true, $false, $false, $false, $true, 1, $false, "79-7=", 2)
$null = $d.Content.Find.Execute("21+17=", $true, $true, $false, $false, $false, $true, 1, $false, "15+68=", 2)
$null = $d.Content.Find.Execute("54+41=", $true, $true, $false, $false, $false, $true, 1, $false, "48+13=", 2)
$null = $d.Content.Find.Execute("83-45=", $true, $true, $false, $false, $false, $true, 1, $false, "18+24=", 2)
$null = $d.Content.Find.Execute("8+61=", $true, $true, $false, $false, $false, $true, 1, $false, "21+2=", 2)
$null = $d.Content.Find.Execute("66-24=", $true, $true, $false, $false, $false, $true, 1, $false, "33+35=", 2)
$null = $d.Content.Find.Execute("34-8=", $true, $true, $false, $false, $false, $true, 1, $false, "35+50=", 2)
$null = $d.Content.Find.Execute("80+5=", $true, $true, $false, $false, $false, $true, 1, $false, "38+35=", 2)
$null = $d.Content.Find.Execute("31+27=", $true, $true, $false, $false, $false, $true, 1, $false, "23+7=", 2)
$null = $d.Content.Find.Execute("22+56=", $true, $true, $false, $false, $false, $true, 1, $false, "97-88=", 2)
$null = $d.Content.Find.Execute("58-19=", $true, $true, $false, $false, $false, $true, 1, $false, "36+52=", 2)
$null = $d.Content.Find.Execute("99-85=", $true, $true, $false, $false, $false, $true, 1, $false, "21+6=", 2)
$null = $d.Content.Find.Execute("32+39=", $true, $true, $false, $false, $false, $true, 1, $false, "17+30=", 2)
$null = $d.Content.Find.Execute("5+58=", $true, $true, $false, $false, $false, $true, 1, $false, "95-87=", 2)
$null = $d.Content.Find.Execute("52-45=", $true, $true, $false, $false, $false, $true, 1, $false, "69-35=", 2)
$null = $d.Content.Find.Execute("48+23=", $true, $true, $false, $false, $false, $true, 1, $false, "72+20=", 2)
$null = $d.Content.Find.Execute("37+3=", $true, $true, $false, $false, $false, $true, 1, $false, "26-17=", 2)
$null = $d.Content.Find.Execute("9+82=", $true, $true, $false, $false, $false, $true, 1, $false, "23+62=", 2)
$null = $d.Content.Find.Execute("9+39=", $true, $true, $false, $false, $false, $true, 1, $false, "75+20=", 2)
$null = $d.Content.Find.Execute("77-44=", $true, $true, $false, $false, $false, $true, 1, $false, "85-28=", 2)
$null = $d.Content.Find.Execute("42+46=", $true, $true, $false, $false, $false, $true, 1, $false, "94-4=", 2)
$null = $d.Content.Find.Execute("29+22=", $true, $true, $false, $false, $false, $true, 1, $false, "49-22=", 2)
$null = $d.Content.Find.Execute("14+64=", $true, $true, $false, $false, $false, $true, 1, $false, "23-14=", 2)
$null = $d.Content.Find.Execute("48+26=", $true, $true, $false, $false, $false, $true, 1, $false, "88-15=", 2)
$null = $d.Content.Find.Execute("34+42=", $true, $true, $false, $false, $false, $true, 1, $false, "99-39=", 2)
$null = $d.Content.Find.Execute("21+70=", $true, $true, $false, $false, $false, $true, 1, $false, "42+50=", 2)
$null = $d.Content.Find.Execute("64-27=", $true, $true, $false, $false, $false, $true, 1, $false, "6+22=", 2)
$null = $d.Content.Find.Execute("21-2=", $true, $true, $false, $false, $false, $true, 1, $false, "45+30=", 2)
$null = $d.Content.Find.Execute("10+34=", $true, $true, $false, $false, $false, $true, 1, $false, "43+42=", 2)
$null = $d.Content.Find.Execute("51-13=", $true, $true, $false, $false, $false, $true, 1, $false, "99-13=", 2)
$null = $d.Content.Find.Execute("47-13=", $true, $true, $false, $false, $false, $true, 1, $false, "54-48=", 2)
$null = $d.Content.Find.Execute("13-8=", $true, $true, $false, $false, $false, $true, 1, $false, "87-69=", 2)
$null = $d.Content.Find.Execute("33+13=", $true, $true, $false, $false, $false, $true, 1, $false, "64-38=", 2)
$null = $d.Content.Find.Execute("52-34=", $true, $true, $false, $false, $false, $true, 1, $false, "93-18=", 2)
$null = $d.Content.Find.Execute("94+5=", $true, $true, $false, $false, $false, $true, 1, $false, "10+59=", 2)
$null = $d.Content.Find.Execute("6+3=", $true, $true, $false, $false, $false, $true, 1, $false, "71-16=", 2)
$null = $d.Content.Find.Execute("6+83=", $true, $true, $false, $false, $false, $true, 1, $false, "64+19=", 2)
$null = $d.Content.Find.Execute("39-18=", $true, $true, $false, $false, $false, $true, 1, $false, "83-60=", 2)
$null = $d.Content.Find.Execute("79-15=", $true, $true, $false, $false, $false, $true, 1, $false, "75-45=", 2)
$null = $d.Content.Find.Execute("7+0=", $true, $true, $false, $false, $false, $true, 1, $false, "37-33=", 2)
$null = $d.Content.Find.Execute("36-7=", $true, $true, $false, $false, $false, $true, 1, $false, "8-7=", 2)
$null = $d.Content.Find.Execute("16+37=", $true, $true, $false, $false, $false, $true, 1, $false, "85-65=", 2)
$null = $d.Content.Find.Execute("6-5=", $true, $true, $false, $false, $false, $true, 1, $false, "95+1=", 2)
$null = $d.Content.Find.Execute("42+16=", $true, $true, $false, $false, $false, $true, 1, $false, "45+4=", 2)
$null = $d.Content.Find.Execute("72-66=", $true, $true, $false, $false, $false, $true, 1, $false, "49+50=", 2)
$null = $d.Content.Find.Execute("66-37=", $true, $true, $false, $false, $false, $true, 1, $false, "11-4=", 2)
$null = $d.Content.Find.Execute("42+40=", $true, $true, $false, $false, $false, $true, 1, $false, "16+6=", 2)
$null = $d.Content.Find.Execute("14+34=", $true, $true, $false, $false, $false, $true, 1, $false, "91-34=", 2)
$null = $d.Content.Find.Execute("95-42=", $true, $true, $false, $false, $false, $true, 1, $false, "79-55=", 2)
$null = $d.Content.Find.Execute("49-28=", $true, $true, $false, $false, $false, $true, 1, $false, "43-18=", 2)
$null = $d.Content.Find.Execute("84+1=", $true, $true, $false, $false, $false, $true, 1, $false, "94-90=", 2)
$null = $d.Content.Find.Execute("80-18=", $true, $true, $false, $false, $false, $true, 1, $false, "26-8=", 2)
$null = $d.Content.Find.Execute("92-74=", $true, $true, $false, $false, $false, $true, 1, $false, "30+34=", 2)
$null = $d.Content.Find.Execute("7+19=", $true, $true, $false, $false, $false, $true, 1, $false, "2+62=", 2)
$null = $d.Content.Find.Execute("51+20=", $true, $true, $false, $false, $false, $true, 1, $false, "60+37=", 2)
$null = $d.Content.Find.Execute("68-30=", $true, $true, $false, $false, $false, $true, 1, $false, "94-21=", 2)
$null = $d.Content.Find.Execute("84-72=", $true, $true, $false, $false, $false, $true, 1, $false, "15+6=", 2)
$null = $d.Content.Find.Execute("41+1=", $true, $true, $false, $false, $false, $true, 1, $false, "32+37=", 2)
$null = $d.Content.Find.Execute("34-31=", $true, $true, $false, $false, $false, $true, 1, $false, "64-0=", 2)
$null = $d.Content.Find.Execute("3+15=", $true, $true, $false, $false, $false, $true, 1, $false, "23+27=", 2)
$null = $d.Content.Find.Execute("24+56=", $true, $true, $false, $false, $false, $true, 1, $false, "41+16=", 2)
$null = $d.Content.Find.Execute("59-24=", $true, $true, $false, $false, $false, $true, 1, $false, "84-45=", 2)
$null = $d.Content.Find.Execute("78-32=", $true, $true, $false, $false, $false, $true, 1, $false, "37+11=", 2)
$null = $d.Content.Find.Execute("91-71=", $true, $true, $false, $false, $false, $true, 1, $false, "94-63=", 2)
$null = $d.Content.Find.Execute("28+46=", $true, $true, $false, $false, $false, $true, 1, $false, "64-56=", 2)
$null = $d.Content.Find.Execute("33-27=", $true, $true, $false, $false, $false, $true, 1, $false, "53+19=", 2)
$null = $d.Content.Find.Execute("54+37=", $true, $true, $false, $false, $false, $true, 1, $false, "68+17=", 2)
$null = $d.Content.Find.Execute("50-2=", $true, $true, $false, $false, $false, $true, 1, $false, "95-67=", 2)
$null = $d.Content.Find.Execute("90-44=", $true, $true, $false, $false, $false, $true, 1, $false, "64-16=", 2)
$null = $d.Content.Find.Execute("59+11=", $true, $true, $false, $false, $false, $true, 1, $false, "78-72=", 2)
$null = $d.Content.Find.Execute("56+36=", $true, $true, $false, $false, $false, $true, 1, $false, "76+3=", 2)
$null = $d.Content.Find.Execute("43+17=", $true, $true, $false, $false, $false, $true, 1, $false, "55-31=", 2)
$null = $d.Content.Find.Execute("50-9=", $true, $true, $false, $false, $false, $true, 1, $false, "23+58=", 2)
$null = $d.Content.Find.Execute("73+19=", $true, $true, $false, $false, $false, $true, 1, $false, "49-7=", 2)
$null = $d.Content.Find.Execute("11+72=", $true, $true, $false, $false, $false, $true, 1, $false, "48-39=", 2)
$null = $d.Content.Find.Execute("24+48=", $true, $true, $false, $false, $false, $true, 1, $false, "6+4=", 2)
$null = $d.Content.Find.Execute("62-61=", $true, $true, $false, $false, $false, $true, 1, $false, "99-84=", 2)
$null = $d.Content.Find.Execute("34+57=", $true, $true, $false, $false, $false, $true, 1, $false, "29+27=", 2)
$null = $d.Content.Find.Execute("18+73=", $true, $true, $false, $false, $false, $true, 1, $false, "32+36=", 2)
$null = $d.Content.Find.Execute("90-76=", $true, $true, $false, $false, $false, $true, 1, $false, "84-6=", 2)
$null = $d.Content.Find.Execute("80-33=", $true, $true, $false, $false, $false, $true, 1, $false, "49+6=", 2)
$null = $d.Content.Find.Execute("98-54=", $true, $true, $false, $false, $false, $true, 1, $false, "11+69=", 2)
$null = $d.Content.Find.Execute("7+7=", $true, $true, $false, $false, $false, $true, 1, $false, "97-46=", 2)
